# Daily attendance processing - swap "Recorded By" name order for System-recorded rows.
# Cells in column G formatted as "System, <email>" are flipped to "<email>, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "System, *") {
        $email = $val.Substring(8)
        $cell.Value2 = $email + ", System"
    }
}
